$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("L1:M1").EntireColumn.Insert()
Write-Output "done"
